# Auto-generated edit script: updates Atomos Profit sheet market-price data
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW/WVR worksheets (scheduled data refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H62").Value = 3860
$ws.Range("I62").Value = 3651.6667
$ws.Range("J62").Value = 4016.25
$ws.Range("K62").Value = 3651.6667
$ws.Range("L62").Value = 4016.25
$ws.Range("M62").Value = -3027.6667
$ws.Range("N62").Value = -5264.25
$ws.Range("H65").Value = 3860
$ws.Range("I65").Value = 3651.6667
$ws.Range("J65").Value = 4016.25
$ws.Range("K65").Value = 18258.3335
$ws.Range("L65").Value = 20081.25
$ws.Range("M65").Value = -15138.3335
$ws.Range("N65").Value = -26321.25
$ws.Range("H86").Value = 2401.4285
$ws.Range("I86").Value = 1761.2
$ws.Range("J86").Value = 4002
$ws.Range("K86").Value = 1761.2
$ws.Range("L86").Value = 4002
$ws.Range("M86").Value = -638.2
$ws.Range("N86").Value = -6248
$ws.Range("H89").Value = 2401.4285
$ws.Range("I89").Value = 1761.2
$ws.Range("J89").Value = 4002
$ws.Range("K89").Value = 8806
$ws.Range("L89").Value = 20010
$ws.Range("M89").Value = -3190
$ws.Range("N89").Value = -31242
$ws.Range("H112").Value = 3969738.2
$ws.Range("J112").Value = 4168008.5
$ws.Range("L112").Value = 12504025.5
$ws.Range("N112").Value = -12506241.5
$ws.Range("H132").Value = 15391685
$ws.Range("I132").Value = 22229102
$ws.Range("K132").Value = 66687306
$ws.Range("M132").Value = -66684776
$ws.Range("H137").Value = 3708367.8
$ws.Range("I137").Value = 4352301
$ws.Range("J137").Value = 5750
$ws.Range("K137").Value = 13056903
$ws.Range("L137").Value = 17250
$ws.Range("M137").Value = -13054353
$ws.Range("N137").Value = -22350
$ws.Range("H141").Value = 563750.9399999999
$ws.Range("I141").Value = 2190.4614
$ws.Range("J141").Value = 928765.25
$ws.Range("K141").Value = 6571.3842
$ws.Range("L141").Value = 2786295.75
$ws.Range("M141").Value = -1391.3842
$ws.Range("N141").Value = -2796655.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4912.81
$ws.Range("I32").Value = 4279.4893
$ws.Range("J32").Value = 14834.833
$ws.Range("K32").Value = 4279.4893
$ws.Range("L32").Value = 14834.833
$ws.Range("M32").Value = -3992.4893
$ws.Range("N32").Value = -15408.833
$ws.Range("H63").Value = 2069.72
$ws.Range("I63").Value = 1511.7142
$ws.Range("J63").Value = 4999.25
$ws.Range("K63").Value = 1511.7142
$ws.Range("L63").Value = 4999.25
$ws.Range("M63").Value = -825.7141999999999
$ws.Range("N63").Value = -6371.25
$ws.Range("H66").Value = 2069.72
$ws.Range("I66").Value = 1511.7142
$ws.Range("J66").Value = 4999.25
$ws.Range("K66").Value = 7558.571
$ws.Range("L66").Value = 24996.25
$ws.Range("M66").Value = -4126.571
$ws.Range("N66").Value = -31860.25
$ws.Range("H74").Value = 2797.4
$ws.Range("I74").Value = 2797.4
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 2797.4
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -1923.4
$ws.Range("N74").ClearContents()
$ws.Range("H77").Value = 2797.4
$ws.Range("I77").Value = 2797.4
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 13987
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -9619
$ws.Range("N77").ClearContents()
$ws.Range("H122").Value = 2080.6667
$ws.Range("I122").Value = 1760.9474
$ws.Range("J122").Value = 3295.6
$ws.Range("K122").Value = 5282.8422
$ws.Range("L122").Value = 9886.799999999999
$ws.Range("M122").Value = -2832.8422
$ws.Range("N122").Value = -14786.8
$ws.Range("H132").Value = 3668.3872
$ws.Range("I132").Value = 3373.3845
$ws.Range("J132").Value = 5202.4
$ws.Range("K132").Value = 10120.1535
$ws.Range("L132").Value = 15607.2
$ws.Range("M132").Value = -7590.1535
$ws.Range("N132").Value = -20667.2

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2098.682
$ws.Range("I105").Value = 1735
$ws.Range("J105").Value = 3068.5
$ws.Range("K105").Value = 1735
$ws.Range("L105").Value = 3068.5
$ws.Range("M105").Value = 12
$ws.Range("N105").Value = -6562.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H86").Value = 10882.5
$ws.Range("I86").Value = 9015.25
$ws.Range("K86").Value = 9015.25
$ws.Range("M86").Value = -7892.25
$ws.Range("H89").Value = 10882.5
$ws.Range("I89").Value = 9015.25
$ws.Range("K89").Value = 45076.25
$ws.Range("M89").Value = -39460.25
$ws.Range("H134").Value = 9333.333000000001
$ws.Range("I134").Value = 0
$ws.Range("J134").Value = 9333.333000000001
$ws.Range("K134").Value = 0
$ws.Range("L134").Value = 27999.999
$ws.Range("M134").ClearContents()
$ws.Range("N134").Value = -33069.999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 1253.7441
$ws.Range("I131").Value = 1361.6666
$ws.Range("J131").Value = 1211.9678
$ws.Range("K131").Value = 4084.9998
$ws.Range("L131").Value = 3635.9034
$ws.Range("M131").Value = 955.0001999999999
$ws.Range("N131").Value = -13715.9034
$ws.Range("H137").Value = 3138.182
$ws.Range("J137").Value = 3440.7693
$ws.Range("L137").Value = 10322.3079
$ws.Range("N137").Value = -20522.3079

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H93").Value = 0
$ws.Range("J93").Value = 0
$ws.Range("L93").Value = 0
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 5932.4
$ws.Range("I132").Value = 10019.429
$ws.Range("J132").Value = 3731.6924
$ws.Range("K132").Value = 30058.287
$ws.Range("L132").Value = 11195.0772
$ws.Range("M132").Value = -27528.287
$ws.Range("N132").Value = -16255.0772
$ws.Range("H138").Value = 29833.334
$ws.Range("J138").Value = 29833.334
$ws.Range("L138").Value = 29833.334
$ws.Range("N138").Value = -40113.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 6462.5454
$ws.Range("I40").Value = 9014.666999999999
$ws.Range("J40").Value = 3400
$ws.Range("K40").Value = 9014.666999999999
$ws.Range("L40").Value = 3400
$ws.Range("M40").Value = -8878.666999999999
$ws.Range("N40").Value = -3672
$ws.Range("H122").Value = 3122.6875
$ws.Range("I122").Value = 2420.3076
$ws.Range("J122").Value = 6166.3335
$ws.Range("K122").Value = 7260.9228
$ws.Range("L122").Value = 18499.0005
$ws.Range("M122").Value = -4810.9228
$ws.Range("N122").Value = -23399.0005
$ws.Range("H141").Value = 29290
$ws.Range("J141").Value = 29290
$ws.Range("L141").Value = 29290
$ws.Range("N141").Value = -39650

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H47").Value = 15000
$ws.Range("J47").Value = 15000
$ws.Range("L47").Value = 15000
$ws.Range("N47").Value = -16144
$ws.Range("H126").Value = 2382540.5
$ws.Range("I126").Value = 1118.037
$ws.Range("J126").Value = 6669101
$ws.Range("K126").Value = 3354.111
$ws.Range("L126").Value = 20007303
$ws.Range("M126").Value = -884.1109999999999
$ws.Range("N126").Value = -20012243
$ws.Range("H132").Value = 207932.86
$ws.Range("I132").Value = 304042.94
$ws.Range("J132").Value = 9705.8125
$ws.Range("K132").Value = 912128.8200000001
$ws.Range("L132").Value = 29117.4375
$ws.Range("M132").Value = -909598.8200000001
$ws.Range("N132").Value = -34177.4375
$ws.Range("H136").Value = 1042.0952
$ws.Range("I136").Value = 507.51514
$ws.Range("J136").Value = 3002.2222
$ws.Range("K136").Value = 1522.54542
$ws.Range("L136").Value = 9006.6666
$ws.Range("M136").Value = 1027.45458
$ws.Range("N136").Value = -14106.6666
$ws.Range("H140").Value = 42187.5
$ws.Range("J140").Value = 42187.5
$ws.Range("L140").Value = 42187.5
$ws.Range("N140").Value = -52547.5

